$wb = $excel.ActiveWorkbook

# Add the two new worksheets at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSame = $wb.Worksheets.Add($null, $lastSheet)
$wsSame.Name = "same_elements"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPartly = $wb.Worksheets.Add($null, $lastSheet)
$wsPartly.Name = "partly_same"

function Fill-Sheet($ws, $rows) {
    $ws.Range("B1").Value = 5
    $ws.Range("C1").Value = 50
    $ws.Range("D1").Value = 500
    $ws.Range("E1").Value = 5000
    $ws.Range("F1").Value = 50000
    $ws.Range("G1").Value = 500000

    $labels = @("byte", "int", "string", "date")
    for ($i = 0; $i -lt 4; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, 1).Value = $labels[$i]
        $rowvals = $rows[$i]
        $ws.Cells.Item($r, 2).Value = $rowvals[0]
        $ws.Cells.Item($r, 3).Value = $rowvals[1]
        $ws.Cells.Item($r, 4).Value = $rowvals[2]
        $ws.Cells.Item($r, 5).Value = $rowvals[3]
        $ws.Cells.Item($r, 6).Value = $rowvals[4]
        $ws.Cells.Item($r, 7).Value = $rowvals[5]
    }
}

$sameRows = @(
    @(0, 0, 0, 0.0005, 0.0005, 0.007507),
    @(0, 0, 0, 0, 0.001001, 0.00901),
    @(0, 0, 0, 0, 0.001001, 0.007508),
    @(0, 0, 0, 0.000501, 0.002003, 0.010511)
)
Fill-Sheet $wsSame $sameRows

$partlyRows = @(
    @(0, 0, 0, 0, 0.001001, 0.008007),
    @(0, 0, 0, 0.000501, 0.001501, 0.029031),
    @(0, 0, 0, 0.0005, 0.004005, 0.04705),
    @(0, 0, 0, 0, 0.002502, 0.041074)
)
Fill-Sheet $wsPartly $partlyRows

# Select the same_elements tab as the active tab
$wsSame.Select()
